$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the "Switches" row (formulas for B/C, plain value for D)
$ws.Range("A2").Value = "Switches"
$ws.Range("B2").Formula = "=4.34+0.1"
$ws.Range("C2").Formula = "=5.235433071+0.1"
$ws.Range("D2").Value = 1.508

# Row 3 becomes "Pot row 1" (values unchanged)
$ws.Range("A3").Value = "Pot row 1"

# Row 4 becomes "Pot row 2" (values unchanged)
$ws.Range("A4").Value = "Pot row 2"

# Row 5 becomes "Pot row 3" with plain values (no formulas)
$ws.Range("A5").Value = "Pot row 3"
$ws.Range("B5").Value = 4.2413385830000001
$ws.Range("C5").Value = 5.137007874
$ws.Range("D5").Value = 4.1500000000000004

# Update the selected cell, matching the saved view state
$ws.Range("C10").Select() | Out-Null
